$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh - apply cell-by-cell updates matching the source diff.
# Column D ("Price") values are textual (dot-grouped / locale-formatted numbers),
# so force Text number format first to stop Excel from re-parsing them as numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.986.87"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.983.76"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.06"
$ws.Range("E5").Value = "  +11.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.92"
$ws.Range("E6").Value = "  +7.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.682"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.78"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.98"
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.621.74"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.999.15"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("E16").Value = "  +8.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.08"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.31"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.642.08"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.01"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.77"
$ws.Range("E22").Value = "  +13.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.00"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  -4.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.22"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.37"
$ws.Range("E26").Value = "  +16.48%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.95"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.31"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.79"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.74"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.03"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "664.80"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.74"
$ws.Range("E36").Value = "  +8.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0899"
$ws.Range("E37").Value = "  +10.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("E42").Value = "  +4.62%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0491"
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.58"
$ws.Range("E45").Value = "  +7.64%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.876.13"
$ws.Range("E49").Value = "  +9.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.02"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("E51").Value = "  +4.19%  "
